$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "25.4."
$ws.Range("C5").Value = "Tomas"
$ws.Range("D5").Value = "Recommender systems reading"

$ws.Range("B6").Select()
